# Auto-generated edit script applying scheduled market-data refresh
# to the Kujata_Profits workbook, per sheet (profession) and row (leve).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2564.1868
$ws.Range("I15").Value = 2564.1868
$ws.Range("K15").Value = 7692.5604
$ws.Range("M15").Value = -7523.5604
$ws.Range("H129").Value = 832.10254
$ws.Range("I129").Value = 550.8
$ws.Range("J129").Value = 873.4706
$ws.Range("K129").Value = 1652.4
$ws.Range("L129").Value = 2620.4118
$ws.Range("M129").Value = 3347.6
$ws.Range("N129").Value = -12620.4118
$ws.Range("H135").Value = 45455870
$ws.Range("I135").Value = 480.29413
$ws.Range("J135").Value = 200004200
$ws.Range("K135").Value = 4322.64717
$ws.Range("L135").Value = 1800037800
$ws.Range("M135").Value = -1787.64717
$ws.Range("N135").Value = -1800042870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2850.8071
$ws.Range("I32").Value = 2324.9143
$ws.Range("K32").Value = 2324.9143
$ws.Range("M32").Value = -2037.9143
$ws.Range("H110").Value = 1261.75
$ws.Range("I110").Value = 835.9
$ws.Range("J110").Value = 1971.5
$ws.Range("K110").Value = 835.9
$ws.Range("L110").Value = 1971.5
$ws.Range("M110").Value = 1209.1
$ws.Range("N110").Value = -6061.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 299
$ws.Range("I2").Value = 299
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 299
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -186
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 318.66666
$ws.Range("I5").Value = 203.5
$ws.Range("J5").Value = 549
$ws.Range("K5").Value = 203.5
$ws.Range("L5").Value = 549
$ws.Range("M5").Value = -91.5
$ws.Range("N5").Value = -773
$ws.Range("H8").Value = 2000
$ws.Range("J8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("N8").Value = -2280
$ws.Range("H10").Value = 515.8
$ws.Range("I10").Value = 515.8
$ws.Range("K10").Value = 515.8
$ws.Range("M10").Value = -376.8
$ws.Range("H11").Value = 1500
$ws.Range("J11").Value = 1500
$ws.Range("L11").Value = 1500
$ws.Range("N11").Value = -1780
$ws.Range("H12").Value = 349.5
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1278
$ws.Range("H16").Value = 1188.1364
$ws.Range("I16").Value = 1224.3572
$ws.Range("J16").Value = 1124.75
$ws.Range("K16").Value = 1224.3572
$ws.Range("L16").Value = 1124.75
$ws.Range("M16").Value = -937.3571999999999
$ws.Range("N16").Value = -1698.75
$ws.Range("H113").Value = 1188.1364
$ws.Range("I113").Value = 1224.3572
$ws.Range("J113").Value = 1124.75
$ws.Range("K113").Value = 1224.3572
$ws.Range("L113").Value = 1124.75
$ws.Range("M113").Value = 945.6428000000001
$ws.Range("N113").Value = -5464.75
$ws.Range("H122").Value = 1067.2
$ws.Range("I122").Value = 1007.3333
$ws.Range("J122").Value = 1157
$ws.Range("K122").Value = 3021.9999
$ws.Range("L122").Value = 3471
$ws.Range("M122").Value = -571.9998999999998
$ws.Range("N122").Value = -8371
$ws.Range("H132").Value = 2132.1428
$ws.Range("I132").Value = 2123.75
$ws.Range("J132").Value = 2143.3333
$ws.Range("K132").Value = 6371.25
$ws.Range("L132").Value = 6429.999899999999
$ws.Range("M132").Value = -3841.25
$ws.Range("N132").Value = -11489.9999
$ws.Range("H134").Value = 25002392
$ws.Range("I134").Value = 2731.6667
$ws.Range("K134").Value = 8195.000100000001
$ws.Range("M134").Value = -5660.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 636.7568
$ws.Range("I113").Value = 539.9167
$ws.Range("J113").Value = 683.24
$ws.Range("K113").Value = 1619.7501
$ws.Range("L113").Value = 2049.72
$ws.Range("M113").Value = 550.2499
$ws.Range("N113").Value = -6389.72
$ws.Range("H131").Value = 17544764
$ws.Range("J131").Value = 1118.2632
$ws.Range("L131").Value = 3354.7896
$ws.Range("N131").Value = -13434.7896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 15085
$ws.Range("J62").Value = 15085
$ws.Range("L62").Value = 15085
$ws.Range("N62").Value = -16457
$ws.Range("H65").Value = 15085
$ws.Range("J65").Value = 15085
$ws.Range("L65").Value = 45255
$ws.Range("N65").Value = -52119
$ws.Range("H70").Value = 26474392
$ws.Range("I70").Value = 20837264
$ws.Range("K70").Value = 20837264
$ws.Range("M70").Value = -20836994
$ws.Range("H73").Value = 26474392
$ws.Range("I73").Value = 20837264
$ws.Range("K73").Value = 20837264
$ws.Range("M73").Value = -20836328
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H132").Value = 3006
$ws.Range("I132").Value = 2760.7896
$ws.Range("K132").Value = 8282.3688
$ws.Range("M132").Value = -5752.3688

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1334.25
$ws.Range("I136").Value = 1088.3077
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 3264.9231
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -714.9231
$ws.Range("N136").Value = -12300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 500
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 500
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H100").Value = 330.45456
$ws.Range("I100").Value = 317.66666
$ws.Range("J100").Value = 345.8
$ws.Range("K100").Value = 635.33332
$ws.Range("L100").Value = 691.6
$ws.Range("M100").Value = -94.33331999999996
$ws.Range("N100").Value = -1773.6
$ws.Range("H132").Value = 1887.1538
$ws.Range("I132").Value = 1659.7894
$ws.Range("K132").Value = 4979.3682
$ws.Range("M132").Value = -2449.3682
$ws.Range("H133").Value = 28475
$ws.Range("J133").Value = 28475
$ws.Range("L133").Value = 28475
$ws.Range("N133").Value = -38595
